$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-18 Friday" "2023-08-19 Saturday"

Replace-Text "47×76=3572" "80×63=5040"
Replace-Text "37×79=2923" "78×24=1872"
Replace-Text "14×13=182" "14×23=322"
Replace-Text "61×25=1525" "13×67=871"
Replace-Text "11×19=209" "50×58=2900"
Replace-Text "23×51=1173" "46×99=4554"
Replace-Text "18×86=1548" "29×31=899"
Replace-Text "50×42=2100" "66×59=3894"
Replace-Text "13×22=286" "11×34=374"
Replace-Text "22×52=1144" "20×78=1560"
Replace-Text "67×68=4556" "71×17=1207"
Replace-Text "82×75=6150" "34×57=1938"
Replace-Text "28×56=1568" "79×98=7742"
Replace-Text "92×84=7728" "63×45=2835"
Replace-Text "71×99=7029" "94×20=1880"
Replace-Text "63×94=5922" "57×73=4161"
Replace-Text "38×43=1634" "50×48=2400"
Replace-Text "28×44=1232" "19×57=1083"
Replace-Text "93×45=4185" "29×44=1276"
Replace-Text "49×28=1372" "78×69=5382"
Replace-Text "49×40=1960" "68×40=2720"
Replace-Text "13×27=351" "55×69=3795"
Replace-Text "90×62=5580" "14×77=1078"
Replace-Text "61×66=4026" "11×38=418"
Replace-Text "95×98=9310" "92×28=2576"
